$wb = $excel.ActiveWorkbook

$wsPed = $wb.Worksheets.Item('Pediatric VFC Vaccine ')
$wsAdult = $wb.Worksheets.Item('Adult Vaccine ')
$wsPedFlu = $wb.Worksheets.Item('Pediatric Influenza Vaccine ')
$wsAdultFlu = $wb.Worksheets.Item('Adult Influenza Vaccine ')

# Pediatric VFC Vaccine
$wsPed.Range('A2').Value2 = 'DTaP '
$wsPed.Range('A3').Value2 = 'DTaP '
$wsPed.Range('A4').Value2 = 'DTaP '
$wsPed.Range('A5').Value2 = 'DTaP-IPV '
$wsPed.Range('A6').Value2 = 'DTaP-IPV '
$wsPed.Range('A7').Value2 = 'DTaP-IPV '
$wsPed.Range('A8').Value2 = 'DTaP-Hep B-IPV '
$wsPed.Range('A9').Value2 = 'DTaP-IP-HI '
$wsPed.Range('A10').Value2 = 'e-IPV '
$wsPed.Range('A11').Value2 = 'Hepatitis A Pediatric '
$wsPed.Range('A12').Value2 = 'Hepatitis A Pediatric '
$wsPed.Range('A13').Value2 = 'Hepatitis A Pediatric '
$wsPed.Range('A14').Value2 = 'Hepatitis A-Hepatitis B 18 only '
$wsPed.Range('A15').Value2 = 'Hepatitis B  Pediatric/Adolescent'
$wsPed.Range('A16').Value2 = 'Hepatitis B  Pediatric/Adolescent'
$wsPed.Range('B16').Value2 = 'Recombivax HB'
$wsPed.Range('A17').Value2 = 'Hib '
$wsPed.Range('A18').Value2 = 'Hib '
$wsPed.Range('A19').Value2 = 'Hib '
$wsPed.Range('A20').Value2 = 'HPV - Human Papillomavirus 9-valent '
$wsPed.Range('A21').Value2 = 'MENB - Meningococcal Group B '
$wsPed.Range('A22').Value2 = 'MENB - Meningococcal Group B '
$wsPed.Range('A23').Value2 = 'MENB - Meningococcal Group B '
$wsPed.Range('A24').Value2 = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$wsPed.Range('A25').Value2 = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$wsPed.Range('A26').Value2 = 'Measles, Mumps and Rubella (MMR) '
$wsPed.Range('A27').Value2 = 'MMR/Varicella '
$wsPed.Range('A28').Value2 = 'Pneumococcal 13-valent  (Pediatric)'
$wsPed.Range('A30').Value2 = 'Rotavirus, Live, Oral, Pentavalent '
$wsPed.Range('A31').Value2 = 'Rotavirus, Live, Oral, Pentavalent '
$wsPed.Range('A32').Value2 = 'Rotavirus, Live, Oral, Oral '
$wsPed.Range('A33').Value2 = 'Tetanus and Diphtheria Toxoids '
$wsPed.Range('A34').Value2 = 'Tetanus and Diphtheria Toxoids '
$wsPed.Range('A35').Value2 = 'Tetanus and Diphtheria Toxoids '
$wsPed.Range('A36').Value2 = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$wsPed.Range('A37').Value2 = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$wsPed.Range('A38').Value2 = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$wsPed.Range('A39').Value2 = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$wsPed.Range('A40').Value2 = 'Varicella '

# Adult Vaccine
$wsAdult.Range('A2').Value2 = 'Hepatitis A-Adult '
$wsAdult.Range('A3').Value2 = 'Hepatitis A-Adult '
$wsAdult.Range('A4').Value2 = 'Hepatitis A Adult '
$wsAdult.Range('A5').Value2 = 'Hepatitis A Adult '
$wsAdult.Range('A6').Value2 = 'Hepatitis A-Hepatitis B Adult '
$wsAdult.Range('A7').Value2 = 'Hepatitis B-Adult '
$wsAdult.Range('A8').Value2 = 'Hepatitis B-Adult '
$wsAdult.Range('A9').Value2 = 'HPV-Human Papillomavirus 9 Valent '
$wsAdult.Range('A10').Value2 = 'Measles, Mumps,  Rubella-Adult '
$wsAdult.Range('A11').Value2 = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$wsAdult.Range('A12').Value2 = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$wsAdult.Range('A13').Value2 = 'MENB - Meningococcal Group B '
$wsAdult.Range('A14').Value2 = 'MENB - Meningococcal Group B '
$wsAdult.Range('A15').Value2 = 'MENB - Meningococcal Group B '
$wsAdult.Range('A16').Value2 = 'Pneumococcal 13-valent  (Adult)'
$wsAdult.Range('A19').Value2 = 'Tetanus and Diphtheria Toxoids '
$wsAdult.Range('A20').Value2 = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$wsAdult.Range('A21').Value2 = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$wsAdult.Range('A22').Value2 = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$wsAdult.Range('A23').Value2 = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$wsAdult.Range('A24').Value2 = 'Varicella-Adult '

# Pediatric Influenza Vaccine
$wsPedFlu.Range('A2').Value2 = 'Influenza  (Age 6 months and older)'
$wsPedFlu.Range('B2').Value2 = 'Fluzone Quadrivalent'
$wsPedFlu.Range('A3').Value2 = 'Influenza  (Age 6-35 months)'
$wsPedFlu.Range('B3').Value2 = 'Fluzone Quadrivalent Pediatric dose'
$wsPedFlu.Range('A4').Value2 = 'Influenza  (Age 36 months and older)'
$wsPedFlu.Range('B4').Value2 = 'Fluzone Quadrivalent'
$wsPedFlu.Range('A5').Value2 = 'Influenza  (Age 36 months and older)'
$wsPedFlu.Range('B5').Value2 = 'Fluzone Quadrivalent'
$wsPedFlu.Range('A6').Value2 = 'Influenza  (Age 6 months and older)'
$wsPedFlu.Range('B6').Value2 = 'Fluarix Quadrivalent'
$wsPedFlu.Range('A7').Value2 = 'Influenza  (Age 6 months and older)'
$wsPedFlu.Range('B7').Value2 = 'FluLaval Quadrivalent'
$wsPedFlu.Range('A8').Value2 = 'Influenza  (Age 6 months and older)'
$wsPedFlu.Range('B8').Value2 = 'FluLaval Quadrivalent'
$wsPedFlu.Range('A9').Value2 = 'Influenza  (Age 4 years and older)'
$wsPedFlu.Range('A10').Value2 = 'Influenza  (Age 4 years and older)'
$wsPedFlu.Range('A11').Value2 = 'Influenza  (Age 5 years and older)'
$wsPedFlu.Range('A12').Value2 = 'Influenza  (Age 5 years and older)'

# Adult Influenza Vaccine
$wsAdultFlu.Range('A2').Value2 = 'Influenza  (Age 6 months and older)'
$wsAdultFlu.Range('B2').Value2 = 'Fluzone Quadrivalent'
$wsAdultFlu.Range('A3').Value2 = 'Influenza  (Age 36 months and older)'
$wsAdultFlu.Range('B3').Value2 = 'Fluzone Quadrivalent'
$wsAdultFlu.Range('A4').Value2 = 'Influenza  (Age 36 months and older)'
$wsAdultFlu.Range('B4').Value2 = 'Fluzone Quadrivalent'
$wsAdultFlu.Range('A5').Value2 = 'Influenza  (Age 6 months and older)'
$wsAdultFlu.Range('B5').Value2 = 'Fluarix Quadrivalent'
$wsAdultFlu.Range('A6').Value2 = 'Influenza  (Age 6 months and older)'
$wsAdultFlu.Range('B6').Value2 = 'FluLaval Quadrivalent'
$wsAdultFlu.Range('A7').Value2 = 'Influenza  (Age 6 months and older)'
$wsAdultFlu.Range('B7').Value2 = 'FluLaval Quadrivalent'
$wsAdultFlu.Range('A8').Value2 = 'Influenza  (Age 4 years and older)'
$wsAdultFlu.Range('A9').Value2 = 'Influenza  (Age 4 years and older)'
$wsAdultFlu.Range('A10').Value2 = 'Influenza  (Age 5 years and older)'
$wsAdultFlu.Range('B10').Value2 = 'Afluria Quadrivalent'
$wsAdultFlu.Range('A11').Value2 = 'Influenza  (Age 5 years and older)'
$wsAdultFlu.Range('B11').Value2 = 'Afluria Quadrivalent'
